# cryptos.xlsx refresh -- GitHub Actions price/volume update.
#
# Column D ("Price") and column E ("Volume(1h)") are stored as literal text
# in this workbook (some prices even use "." as a thousands separator, e.g.
# "64.857.68"), and row 48/49 swap which coin (Arweave / InjectiveProtocol)
# occupies which row. Column E values keep their original two-space padding.
#
# Range.Value auto-coerces a plain decimal-looking string (e.g. "616.81")
# into a real number. To keep those cells as text -- matching the original
# inline-string cells -- Set-TextValue briefly switches the cell to the
# "Text" number format, assigns the string, then restores the "Normal" cell
# style so no stray formatting is left behind.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $text) {
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = "Normal"
}

# Row 2 - Bitcoin
$ws.Range("D2").Value = "64.857.68"
$ws.Range("E2").Value = "  +1.18%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "3.173.07"
$ws.Range("E3").Value = "  +1.37%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  -0.04%  "

# Row 5 - BNB
Set-TextValue $ws.Range("D5") "616.81"
$ws.Range("E5").Value = "  +1.31%  "

# Row 6 - Solana
Set-TextValue $ws.Range("D6") "146.86"
$ws.Range("E6").Value = "  -1.33%  "

# Row 7 - USDC
$ws.Range("E7").Value = "  -0.04%  "

# Row 8 - LidoStakedEther
$ws.Range("D8").Value = "3.168.83"
$ws.Range("E8").Value = "  +1.21%  "

# Row 9 - XRP
$ws.Range("E9").Value = "  -0.29%  "

# Row 10 - Dogecoin
$ws.Range("E10").Value = "  +0.33%  "

# Row 11 - Toncoin
Set-TextValue $ws.Range("D11") "5.49"
$ws.Range("E11").Value = "  -0.94%  "

# Row 12 - Cardano
$ws.Range("E12").Value = "  -0.31%  "

# Row 13 - ShibaInu
$ws.Range("E13").Value = "  +1.32%  "

# Row 14 - Avalanche
Set-TextValue $ws.Range("D14") "35.93"
$ws.Range("E14").Value = "  -2.25%  "

# Row 15 - WrappedliquidstakedEther2.0
$ws.Range("D15").Value = "3.696.99"
$ws.Range("E15").Value = "  +1.09%  "

# Row 16 - TRON
$ws.Range("E16").Value = "  +3.11%  "

# Row 17 - WrappedBTC
$ws.Range("D17").Value = "64.868.53"
$ws.Range("E17").Value = "  +1.01%  "

# Row 18 - WrappedEther
$ws.Range("D18").Value = "3.172.67"
$ws.Range("E18").Value = "  +1.23%  "

# Row 19 - Polkadot
Set-TextValue $ws.Range("D19") "6.95"
$ws.Range("E19").Value = "  -0.52%  "

# Row 20 - BitcoinCash
Set-TextValue $ws.Range("D20") "480.65"
$ws.Range("E20").Value = "  -0.36%  "

# Row 21 - Chainlink
Set-TextValue $ws.Range("D21") "14.76"
$ws.Range("E21").Value = "  +1.03%  "

# Row 22 - Polygon
$ws.Range("E22").Value = "  +1.49%  "

# Row 23 - Uniswap
Set-TextValue $ws.Range("D23") "7.99"
$ws.Range("E23").Value = "  +2.98%  "

# Row 24 - InternetComputer(DFINITY)
Set-TextValue $ws.Range("D24") "13.84"
$ws.Range("E24").Value = "  +0.18%  "

# Row 25 - Litecoin
Set-TextValue $ws.Range("D25") "84.67"
$ws.Range("E25").Value = "  +0.84%  "

# Row 26 - Dai
$ws.Range("E26").Value = "  +0.00%  "

# Row 27 - PancakeSwap
Set-TextValue $ws.Range("D27") "2.83"
$ws.Range("E27").Value = "  -3.58%  "

# Row 28 - RenderToken
Set-TextValue $ws.Range("D28") "8.68"
$ws.Range("E28").Value = "  +1.45%  "

# Row 29 - Hedera
$ws.Range("E29").Value = "  -4.92%  "

# Row 30 - ImmutableX
Set-TextValue $ws.Range("D30") "2.11"
$ws.Range("E30").Value = "  -5.62%  "

# Row 31 - NEARProtocol
$ws.Range("E31").Value = "  -0.27%  "

# Row 32 - FirstDigitalUSD
$ws.Range("E32").Value = "  +0.04%  "

# Row 33 - Stacks
$ws.Range("E33").Value = "  +0.43%  "

# Row 34 - EthereumClassic
Set-TextValue $ws.Range("D34") "26.74"
$ws.Range("E34").Value = "  +0.05%  "

# Row 35 - Mantle
$ws.Range("E35").Value = "  +2.55%  "

# Row 36 - PEPE
$ws.Range("E36").Value = "  +5.86%  "

# Row 37 - Filecoin
Set-TextValue $ws.Range("D37") "6.05"
$ws.Range("E37").Value = "  -0.59%  "

# Row 38 - dogwifhat
$ws.Range("E38").Value = "  -0.16%  "

# Row 39 - OKB
Set-TextValue $ws.Range("D39") "53.21"
$ws.Range("E39").Value = "  -2.48%  "

# Row 40 - Bittensor
Set-TextValue $ws.Range("D40") "467.38"
$ws.Range("E40").Value = "  +3.82%  "

# Row 41 - VeChain
$ws.Range("E41").Value = "  +0.69%  "

# Row 42 - Kaspa
$ws.Range("E42").Value = "  -2.82%  "

# Row 43 - Cosmos
Set-TextValue $ws.Range("D43") "8.42"
$ws.Range("E43").Value = "  -0.75%  "

# Row 44 - Maker
$ws.Range("D44").Value = "2.860.11"
$ws.Range("E44").Value = "  -0.53%  "

# Row 45 - Fetch.AI
$ws.Range("E45").Value = "  +0.61%  "

# Row 46 - TheGraph
$ws.Range("E46").Value = "  -0.73%  "

# Row 47 - ThetaToken
Set-TextValue $ws.Range("D47") "2.45"
$ws.Range("E47").Value = "  +6.39%  "

# Row 48 - Arweave / InjectiveProtocol (rows swapped)
$ws.Range("B48").Value = "Arweave"
$ws.Range("C48").Value = "https://coinranking.com/coin/7XWg41D1+arweave-ar"
Set-TextValue $ws.Range("D48") "37.07"
$ws.Range("E48").Value = "  +11.85%  "

# Row 49 - InjectiveProtocol / Arweave (rows swapped)
$ws.Range("B49").Value = "InjectiveProtocol"
$ws.Range("C49").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
Set-TextValue $ws.Range("D49") "26.85"
$ws.Range("E49").Value = "  +0.92%  "

# Row 50 - USDe
$ws.Range("E50").Value = "  +0.12%  "

# Row 51 - Stellar
$ws.Range("E51").Value = "  -0.72%  "
